# Update the EVD linelist with a new version of the data: new id hashes in
# column A, shifted onset dates in column B, re-shuffled sex in column C,
# and new age values in column D, for rows 2-13.
#
# Some of the new id hashes (e.g. "664549", "947e40", "185911", "605322")
# parse as plain numbers/scientific-notation, so a naive .Value assignment
# would have Excel's type-sniffing store them as numbers instead of text,
# which doesn't match how every other "id" cell in the sheet is stored
# (shared-string text, default/unstyled cell). To force text without
# leaving any stray cell formatting behind, each id is written through a
# text-returning formula and then "flattened" to a literal value via
# Copy + PasteSpecial(values only).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-IdValue($range, $value) {
    $esc = $value -replace '"', '""'
    $range.Formula = "=""$esc"""
    $range.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
}

Set-IdValue $ws.Range("A2") "39e9dc"
$ws.Range("B2").Value = 43018.0
$ws.Range("C2").Value = "female"
$ws.Range("D2").Value = 62.0

Set-IdValue $ws.Range("A3") "664549"
$ws.Range("B3").Value = 43024.0
$ws.Range("C3").Value = "male"
$ws.Range("D3").Value = 28.0

Set-IdValue $ws.Range("A4") "b4d8aa"
$ws.Range("B4").Value = 43025.0
$ws.Range("C4").Value = "male"
$ws.Range("D4").Value = 54.0

Set-IdValue $ws.Range("A5") "51883d"
$ws.Range("B5").Value = 43026.0
$ws.Range("C5").Value = "male"
$ws.Range("D5").Value = 57.0

Set-IdValue $ws.Range("A6") "947e40"
$ws.Range("B6").Value = 43028.0
$ws.Range("C6").Value = "female"
$ws.Range("D6").Value = 23.0

Set-IdValue $ws.Range("A7") "9aa197"
$ws.Range("B7").Value = 43028.0
$ws.Range("C7").Value = "female"
$ws.Range("D7").Value = 66.0

Set-IdValue $ws.Range("A8") "e4b0a2"
$ws.Range("B8").Value = 43029.0
$ws.Range("C8").Value = "female"
$ws.Range("D8").Value = 13.0

Set-IdValue $ws.Range("A9") "af0ac0"
$ws.Range("B9").Value = 43029.0
$ws.Range("C9").Value = "male"
$ws.Range("D9").Value = 10.0

Set-IdValue $ws.Range("A10") "185911"
$ws.Range("B10").Value = 43029.0
$ws.Range("C10").Value = "female"
$ws.Range("D10").Value = 34.0

Set-IdValue $ws.Range("A11") "601d2e"
$ws.Range("B11").Value = 43030.0
$ws.Range("C11").Value = "male"
$ws.Range("D11").Value = 11.0

Set-IdValue $ws.Range("A12") "605322"
$ws.Range("B12").Value = 43030.0
$ws.Range("C12").Value = "female"
$ws.Range("D12").Value = 23.0

Set-IdValue $ws.Range("A13") "e399b1"
$ws.Range("B13").Value = 43031.0
$ws.Range("C13").Value = "female"
$ws.Range("D13").Value = 23.0
